{"js": "// Add test for written numbers rule:\n// In the \"BROTTET\" paragraph, turn the full stop after the first\n// \"m\u00e5ls\u00e4gande\" into a comma and insert a new clause\n// \"de kanske var 3 m\u00e4n, i v\u00e4rsta fall 14\" before the following sentence\n// \"Misshandeln har best\u00e5tt av ...\".\nconst body = context.document.body;\n\nconst results = body.search(\"m\u00e5ls\u00e4gande. Misshandeln\", { matchCase: true, matchWholeWord: false });\nresults.load(\"items\");\nawait context.sync();\n\nif (results.items.length === 0) {\n  throw new Error(\"Target text not found\");\n}\n\nresults.items[0].insertText(\n  \"m\u00e5ls\u00e4gande, de kanske var 3 m\u00e4n, i v\u00e4rsta fall 14. Misshandeln\",\n  Word.InsertLocation.replace\n);\n\nawait context.sync();\n", "ps1": "# Add test for written numbers rule:\n# In the \"BROTTET\" paragraph, turn the full stop after the first\n# \"m\u00e5ls\u00e4gande\" into a comma and insert a new clause\n# \"de kanske var 3 m\u00e4n, i v\u00e4rsta fall 14\" before the following sentence\n# \"Misshandeln har best\u00e5tt av ...\".\n\n$d = $word.ActiveDocument\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"m\u00e5ls\u00e4gande. Misshandeln\"\n$find.Replacement.Text = \"m\u00e5ls\u00e4gande, de kanske var 3 m\u00e4n, i v\u00e4rsta fall 14. Misshandeln\"\n$find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2)\n"}
